$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing Text storage (the sheet's
# Price/Volume columns hold text like '216.35' or '  +1.06%  ', and a plain
# .Value assignment would let Excel auto-coerce number-looking strings into
# actual numbers). Apply a Text number format for the write, then restore the
# cell's plain/default style so no stray formatting is introduced.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row-by-row updates to Price (D) and Volume(1h) (E) columns
Set-TextValue "D2" '28.026.26'
Set-TextValue "E2" '  +3.24%  '

Set-TextValue "D3" '1.686.64'
Set-TextValue "E3" '  +0.45%  '

Set-TextValue "E4" '  -0.10%  '

Set-TextValue "D5" '216.22'
Set-TextValue "E5" '  +0.93%  '

Set-TextValue "D6" '0.523'
Set-TextValue "E6" '  +0.94%  '

Set-TextValue "E7" '  -0.13%  '

Set-TextValue "D8" '23.71'
Set-TextValue "E8" '  +4.14%  '

Set-TextValue "E9" '  +1.54%  '

Set-TextValue "E10" '  +0.72%  '

Set-TextValue "D11" '0.0885'
Set-TextValue "E11" '  -0.61%  '

Set-TextValue "D12" '1.927.08'
Set-TextValue "E12" '  +0.55%  '

Set-TextValue "D13" '1.686.73'
Set-TextValue "E13" '  +0.44%  '

Set-TextValue "D14" '4.18'
Set-TextValue "E14" '  -0.18%  '

Set-TextValue "D15" '0.556'
Set-TextValue "E15" '  +1.07%  '

Set-TextValue "D16" '66.96'
Set-TextValue "E16" '  +0.56%  '

Set-TextValue "D17" '250.94'
Set-TextValue "E17" '  +6.54%  '

Set-TextValue "D18" '28.006.74'
Set-TextValue "E18" '  +3.30%  '

Set-TextValue "E19" '  +0.28%  '

Set-TextValue "D20" '7.62'
Set-TextValue "E20" '  -3.38%  '

Set-TextValue "E21" '  -0.14%  '

Set-TextValue "D22" '4.52'
Set-TextValue "E22" '  -0.39%  '

Set-TextValue "D23" '9.54'
Set-TextValue "E23" '  +0.10%  '

Set-TextValue "E24" '  -1.71%  '

Set-TextValue "D25" '147.38'
Set-TextValue "E25" '  +0.31%  '

Set-TextValue "D26" '7.32'
Set-TextValue "E26" '  -1.28%  '

Set-TextValue "D27" '16.46'
Set-TextValue "E27" '  +0.81%  '

Set-TextValue "E28" '  +0.53%  '

Set-TextValue "E29" '  -0.15%  '

Set-TextValue "E30" '  +6.92%  '

Set-TextValue "D31" '0.0503'
Set-TextValue "E31" '  +0.27%  '

Set-TextValue "E32" '  +0.37%  '

Set-TextValue "E33" '  -1.95%  '

Set-TextValue "D34" '1.428.51'
Set-TextValue "E34" '  -7.34%  '

Set-TextValue "D35" '1.59'
Set-TextValue "E35" '  -4.01%  '

Set-TextValue "D36" '0.937'
Set-TextValue "E36" '  -0.40%  '

Set-TextValue "E37" '  -0.31%  '

Set-TextValue "D38" '0.590'
Set-TextValue "E38" '  -2.79%  '

Set-TextValue "E39" '  -0.01%  '

Set-TextValue "E40" '  -3.05%  '

Set-TextValue "D41" '69.71'
Set-TextValue "E41" '  +0.34%  '

Set-TextValue "E42" '  -0.13%  '

Set-TextValue "E43" '  -5.42%  '

Set-TextValue "E46" '  +2.70%  '

Set-TextValue "E47" '  +5.42%  '

Set-TextValue "D48" '89.27'
Set-TextValue "E48" '  -0.59%  '

Set-TextValue "D49" '0.0₆0111'
Set-TextValue "E49" '  -0.74%  '

Set-TextValue "D50" '0.102'
Set-TextValue "E50" '  -1.28%  '

Set-TextValue "D51" '7.83'
Set-TextValue "E51" '  -4.48%  '

# Rows 44 and 45: coin order swapped (RocketPoolETH <-> MXToken) together with
# refreshed price/volume figures
Set-TextValue "B44" 'MXToken'
Set-TextValue "C44" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D44" '2.24'
Set-TextValue "E44" '  -0.66%  '

Set-TextValue "B45" 'RocketPoolETH'
Set-TextValue "C45" 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue "D45" '1.833.55'
Set-TextValue "E45" '  +0.53%  '

